# Adds an "original word count" column (G) to the "backstories" sheet that
# counts the words in column F (Original Backstory) using the classic
# LEN/SUBSTITUTE word-count trick, plus summary statistics (STDEV.S) for
# both the existing word-count column (D) and the new one (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("backstories")

# New header for column G
$ws.Range("G1").Value = "original word count"

# Per-row word-count formulas for column G (mirrors column F's text)
$ws.Range("G2").Formula = '=LEN(TRIM(F2)) - LEN(SUBSTITUTE(F2, " ", "")) + 1'
$ws.Range("G3:G13").Formula = '=LEN(TRIM(F3)) - LEN(SUBSTITUTE(F3, " ", "")) + 1'

# Summary statistics row (row 14)
$ws.Range("E14").Formula = "=  _xlfn.STDEV.S(D2:D13)"
$ws.Range("G14").Formula = " =AVERAGE(G2:G8, G9:G13)"
$ws.Range("H14").Formula = "=  _xlfn.STDEV.S(G2:G13)"

# Column widths for the two new columns, matching the edited workbook
$ws.Columns.Item(7).ColumnWidth = 18.33
$ws.Columns.Item(8).ColumnWidth = 15.78
